# Add "Wins", "Losses", "Ties" columns (AD, AE, AF) to Sheet1.
# Header cells (row 1) get the same formatting as the existing header row
# by copying the adjacent header cell's format; data rows (2-40) get a
# constant season record of 89 wins, 73 losses, 0 ties.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Header row ---------------------------------------------------------
$ws.Range("AC1").Copy($ws.Range("AD1"))
$ws.Range("AC1").Copy($ws.Range("AE1"))
$ws.Range("AC1").Copy($ws.Range("AF1"))

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (2-40): season record -----------------------------------
$wins = 89
$losses = 73
$ties = 0

for ($row = 2; $row -le 40; $row++) {
    $ws.Cells.Item($row, 30).Value = $wins
    $ws.Cells.Item($row, 31).Value = $losses
    $ws.Cells.Item($row, 32).Value = $ties
}

Write-Output "Added Wins/Losses/Ties columns (AD:AF) for rows 1-40"
